# Natmi following Dr Hou advice
# Update existing LR-pair rows (FAPs -> Fgf10/Fgfrl1 -> {ECs, FAPs, sCs})
# with revised numbers, re-label row 4's target cluster from "sCs" to "M2",
# and append a new row 5 for target cluster "sCs" with its own numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: FAPs / Fgf10 / Fgfrl1 / ECs ---
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.31623
$ws.Range("H2").Value = 3.94869
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.704056
$ws.Range("N2").Value = 2.112168
$ws.Range("O2").Value = 0.06182348620901491
$ws.Range("P2").Value = 0.0618234862090149
$ws.Range("Q2").Value = 0.92669962888
$ws.Range("R2").Value = 8.34029665992
$ws.Range("S2").Value = 0.06182348620901491
$ws.Range("T2").Value = 0.0618234862090149

# --- Row 3: FAPs / Fgf10 / Fgfrl1 / FAPs ---
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.31623
$ws.Range("H3").Value = 3.94869
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 8.646516666666667
$ws.Range("N3").Value = 25.93955
$ws.Range("O3").Value = 0.759254667097055
$ws.Range("P3").Value = 0.759254667097055
$ws.Range("Q3").Value = 11.38080463216667
$ws.Range("R3").Value = 102.4272416895
$ws.Range("S3").Value = 0.759254667097055
$ws.Range("T3").Value = 0.759254667097055

# --- Row 4: FAPs / Fgf10 / Fgfrl1 / M2 (re-clustered from sCs) ---
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.31623
$ws.Range("H4").Value = 3.94869
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.4191446666666667
$ws.Range("N4").Value = 1.257434
$ws.Range("O4").Value = 0.03680528895322079
$ws.Range("P4").Value = 0.03680528895322079
$ws.Range("Q4").Value = 0.5516907846066667
$ws.Range("R4").Value = 4.96521706146
$ws.Range("S4").Value = 0.03680528895322079
$ws.Range("T4").Value = 0.03680528895322079

# --- Row 5 (new): FAPs / Fgf10 / Fgfrl1 / sCs ---
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fgf10"
$ws.Range("C5").Value = "Fgfrl1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.31623
$ws.Range("H5").Value = 3.94869
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.618446666666667
$ws.Range("N5").Value = 4.85534
$ws.Range("O5").Value = 0.1421165577407093
$ws.Range("P5").Value = 0.1421165577407092
$ws.Range("Q5").Value = 2.130248056066667
$ws.Range("R5").Value = 19.1722325046
$ws.Range("S5").Value = 0.1421165577407093
$ws.Range("T5").Value = 0.1421165577407092
